$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137, pushing existing rows 137-145 down to 138-146.
$ws.Rows.Item(137).Insert()

# The new row 137 shares most field values with its neighbours (constant
# across this subset); only the date, volume/price and origin columns
# differ for this new weekly record.
$ws.Cells.Item(137, 1).Value = 11
$ws.Cells.Item(137, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(137, 3).Value = "Bíobío"
$ws.Cells.Item(137, 4).Value = 45013
$ws.Cells.Item(137, 5).Value = 8
$ws.Cells.Item(137, 6).Value = 100112001
$ws.Cells.Item(137, 7).Value = "Berenjena"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 220
$ws.Cells.Item(137, 11).Value = 8000
$ws.Cells.Item(137, 12).Value = 8500
$ws.Cells.Item(137, 13).Value = 8273
$ws.Cells.Item(137, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(137, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(137, 16).Value = 138
$ws.Cells.Item(137, 17).Value = 60
$ws.Cells.Item(137, 18).Value = "Hortaliza"
